# Apply updated Betfair Back/Lay odds values for 2025-12-23
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Melbourne City vs Macarthur FC
$ws.Range("H2").Value = 5.7
$ws.Range("I2").Value = 5.8
$ws.Range("O2").Value = 1.27
$ws.Range("R2").Value = 1.44
$ws.Range("T2").Value = 1.83
$ws.Range("U2").Value = 2.1
$ws.Range("V2").Value = 1.2
$ws.Range("AC2").Value = 9.199999999999999
$ws.Range("AG2").Value = 9.6
$ws.Range("AN2").Value = 9.6

# Row 3 - MC Alger vs ES Ben Aknoun
$ws.Range("G3").Value = 1.48
$ws.Range("H3").Value = 11.5
$ws.Range("I3").Value = 15.5
$ws.Range("J3").Value = 4.1
$ws.Range("K3").Value = 4.9
$ws.Range("O3").Value = 1.45
$ws.Range("S3").Value = 4.6
$ws.Range("U3").Value = 1.49
$ws.Range("W3").Value = 3.05
$ws.Range("Y3").Value = 32
$ws.Range("AC3").Value = 11.5
$ws.Range("AD3").Value = 65
$ws.Range("AF3").Value = 7
$ws.Range("AG3").Value = 12
$ws.Range("AN3").Value = 1000

# Row 4 - Guimaraes vs Sporting Lisbon
$ws.Range("L4").Value = 1.4
$ws.Range("V4").Value = 3.3
$ws.Range("AJ4").Value = 380
